$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [B, C, D, E] new values
$data = @{
    2  = @(0.5454545454545454, 0.4615384615384616, 0.4999999999999999, 13)
    3  = @(0.4615384615384616, 0.5454545454545454, 0.4999999999999999, 11)
    4  = @(0.5, 0.5, 0.5, 0.5)
    5  = @(0.5034965034965035, 0.5034965034965035, 0.4999999999999999, 24)
    6  = @(0.506993006993007, 0.5, 0.4999999999999998, 24)
    7  = @(0.4166666666666667, 0.3846153846153846, 0.4, 13)
    8  = @(0.3333333333333333, 0.3636363636363636, 0.3478260869565217, 11)
    9  = @(0.375, 0.375, 0.375, 0.375)
    10 = @(0.375, 0.3741258741258742, 0.3739130434782609, 24)
    11 = @(0.3784722222222223, 0.375, 0.3760869565217391, 24)
    12 = @(0.5384615384615384, 0.5384615384615384, 0.5384615384615384, 13)
    13 = @(0.4545454545454545, 0.4545454545454545, 0.4545454545454545, 11)
    14 = @(0.5, 0.5, 0.5, 0.5)
    15 = @(0.4965034965034965, 0.4965034965034965, 0.4965034965034965, 24)
    16 = @(0.5, 0.5, 0.5, 24)
    17 = @(0.6153846153846154, 0.6153846153846154, 0.6153846153846154, 13)
    18 = @(0.5454545454545454, 0.5454545454545454, 0.5454545454545454, 11)
    20 = @(0.5804195804195804, 0.5804195804195804, 0.5804195804195804, 24)
    21 = @(0.5833333333333334, 0.5833333333333334, 0.5833333333333334, 24)
    22 = @(0.5882352941176471, 0.7692307692307693, 0.6666666666666667, 13)
    23 = @(0.5714285714285714, 0.3636363636363636, 0.4444444444444444, 11)
    24 = @(0.5833333333333334, 0.5833333333333334, 0.5833333333333334, 0.5833333333333334)
    25 = @(0.5798319327731092, 0.5664335664335665, 0.5555555555555556, 24)
    26 = @(0.5805322128851541, 0.5833333333333334, 0.5648148148148149, 24)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
